$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 154107
$ws.Range("C5").Value = 8787
$ws.Range("C6").Value = 846
$ws.Range("C7").Value = 5.7
